$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(8,5,4,1,7,5,3,7,9,10,8,3,7,6,8,7,6,8,8,1,6,10,9,7,10,6,1,1,4,5,3,6,8,10,8,9,2,10,5,5,2,5,1,7,9,4,9,1,2,1,8,2,1,5,1,6,10,10,6,6,1,5,4,7,2,9,5,7,8,10,3,1,3,5,9,8,10,2,7,9,6,5,2,8,6,6,4,4,10,3,5,1,3,6,7,3,4,2,6,2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

for ($r = 101; $r -le 110; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=IF(A" + $r + ">=10,TRUE())*1"
    $ws.Cells.Item($r, 4).Formula = "=IF(C" + $r + "=1,B" + $r + ")*1"
}
